# Remove the "H 72" record (originally row 2) from the data sheet.
# Excel shifts the remaining rows (3:63) up by one automatically, which
# also updates the sheet's used-range dimension from A1:F63 to A1:F62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
